$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'291.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.23%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'31.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.81%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.958"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.26%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07453"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.74%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.251"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-5.63%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.740"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.96%"
$ws.Range("E7").Style = "Normal"

$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.763"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.57%"
$ws.Range("E8").Style = "Normal"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9202"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.63%"
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09375"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'17.89%"
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1726"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.49%"
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.70%"
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03219"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'4.20%"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09930"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.78%"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005726"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.94%"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.477"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.18%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.130"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.63%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'0.33%"
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'0.25%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.178"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'5.39%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.2120"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-7.87%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04514"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.10%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001218"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.66%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004258"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-3.40%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001299"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.05%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0003391"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.08%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").Value = "'0.01621"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.08%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04576"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.25%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007436"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.38%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.009828"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = "'0.1358"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.45%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002157"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'6.70%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.008724"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-8.20%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006102"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.54%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'2.525"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'12.69%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.001998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-31.06%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002098"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0001998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.08%"
$ws.Range("E51").Style = "Normal"
